# week 2 slides and r code
#
# Slide 21 ("Other important terms") - "Content Placeholder 2":
# bold the lvl=2 bullet "Role of experimentation vs observation"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(21)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The 6th paragraph in the placeholder is the target bullet.
$target = $tr.Paragraphs(6, 1)
$target.Font.Bold = $true
